$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 28) - a new restaurant entry
$ws.Range("A28").Value = "Neuse River Brewing"
$ws.Range("B28").Value = "Raleigh"
$ws.Range("C28").Value = "Neuse River Burger"
$ws.Range("D28").Value = "Brassiere/Burgers"
$ws.Range("E28").Value = 35.80457
$ws.Range("F28").Value = -78.6325

# Update selection / view to reflect the new last row
$ws.Range("A28:F28").Select()
